$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New row 11: "Pliki - zapis do tabel" ---
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Pliki - zapis do tabel"
$ws.Range("C11").Value = "doprecyzowanie pól do zapisu"
$ws.Range("D11").Value = "Wordy\Pliki zapisy do tabel.docx"

# Match the wrap-text formatting used by the other "opis" cells (col C)
$ws.Range("C11").WrapText = $true

# Row height to match the other two-line rows
$ws.Rows.Item(11).RowHeight = 28.5

# Hyperlink for the new docx reference, then re-apply the hyperlink style
# (matches the styling used by the other hyperlinked cells in column D)
$ws.Hyperlinks.Add($ws.Range("D11"), "Wordy\Pliki zapisy do tabel.docx")
$ws.Range("D11").Style = $ws.Range("D10").Style

# Update the selection to the new hyperlink cell and bring it into view
$null = $ws.Range("D11").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
